$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.377.63'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '2.750.35'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '''574.56'
$ws.Range('E5').Value = '  -0.91%  '
$ws.Range('D6').Value = '''159.15'
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '''0.601'
$ws.Range('E8').Value = '  -1.94%  '
$ws.Range('E9').Value = '  -2.07%  '
$ws.Range('E10').Value = '  +4.96%  '
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '''0.386'
$ws.Range('E12').Value = '  -1.84%  '
$ws.Range('D13').Value = '3.240.12'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').Value = '''27.15'
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('D15').Value = '64.050.78'
$ws.Range('E15').Value = '  +0.30%  '
$ws.Range('E16').Value = '  -2.47%  '
$ws.Range('D17').Value = '2.757.80'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').Value = '''12.14'
$ws.Range('E18').Value = '  -1.44%  '
$ws.Range('D19').Value = '''4.83'
$ws.Range('E19').Value = '  -2.83%  '
$ws.Range('D20').Value = '''356.99'
$ws.Range('E20').Value = '  -1.62%  '
$ws.Range('D21').Value = '''6.65'
$ws.Range('E21').Value = '  -3.65%  '
$ws.Range('D22').Value = '''0.998'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').Value = '''0.526'
$ws.Range('E23').Value = '  -8.06%  '
$ws.Range('D24').Value = '''64.91'
$ws.Range('E24').Value = '  -2.26%  '
$ws.Range('E25').Value = '  -1.09%  '
$ws.Range('D26').Value = '''8.56'
$ws.Range('E26').Value = '  -1.45%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').Value = '0.0₃0919'
$ws.Range('E28').Value = '  -2.50%  '
$ws.Range('D29').Value = '''7.34'
$ws.Range('E29').Value = '  +3.13%  '
$ws.Range('D30').Value = '''1.37'
$ws.Range('E30').Value = '  +8.63%  '
$ws.Range('E31').Value = '  -1.87%  '
$ws.Range('D32').Value = '''167.37'
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('D33').Value = '''4.97'
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('D34').Value = '''1.50'
$ws.Range('E34').Value = '  +1.83%  '
$ws.Range('D35').Value = '''20.16'
$ws.Range('E35').Value = '  -2.11%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').Value = '''1.83'
$ws.Range('E37').Value = '  +0.75%  '
$ws.Range('D38').Value = '''0.998'
$ws.Range('E38').Value = '  -1.98%  '
$ws.Range('D39').Value = '''350.05'
$ws.Range('E39').Value = '  +4.44%  '
$ws.Range('E40').Value = '  +4.02%  '
$ws.Range('D41').Value = '''4.17'
$ws.Range('E41').Value = '  -1.39%  '
$ws.Range('D42').Value = '''39.09'
$ws.Range('E42').Value = '  -1.14%  '
$ws.Range('D43').Value = '''22.48'
$ws.Range('E43').Value = '  +1.95%  '
$ws.Range('D44').Value = '''21.53'
$ws.Range('E44').Value = '  -2.76%  '
$ws.Range('D45').Value = '''0.0591'
$ws.Range('E45').Value = '  -1.65%  '
$ws.Range('D46').Value = '''135.98'
$ws.Range('E46').Value = '  -0.57%  '
$ws.Range('D47').Value = '''0.628'
$ws.Range('E47').Value = '  -1.99%  '
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('D49').Value = '''0.0252'
$ws.Range('E49').Value = '  -2.66%  '
$ws.Range('E50').Value = '  -0.06%  '
$ws.Range('D51').Value = '2.138.98'
$ws.Range('E51').Value = '  +0.27%  '
